# Rename the worksheet from "Sheet1" to "sheet"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "sheet"

# Re-label some "iphone" rows as "samsung" to create a per-product
# dictionary split (adds a new shared string "samsung").
$ws.Range("A11").Value = "samsung"
$ws.Range("A19").Value = "samsung"

# Update the active cell/selection to A12 (was B20).
$ws.Range("A12").Select()
